# EEBEGU-720 - adding superscript to footnotes in Finanzielle Situation Dokument
#
# Replace the bracketed footnote markers "[1]" / "[2]" used throughout the
# "Berechnungsgrundlagen" tables (and their footer legend) with the actual
# Unicode superscript glyphs "¹" / "²", and relocate the stray _GoBack
# bookmark to the position of the very last edit (end of the last
# "Zwischentotal Nettovermögen insgesamt" cell).

$d = $word.ActiveDocument

# --- 1. Body: " [1]" -> " ¹" --------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Execute(" [1]", $false, $false, $false, $false, $false, $true, 1, $false, " ¹", 2)

# --- 2. Body: "[2]" -> " ²" -----------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Execute("[2]", $false, $false, $false, $false, $false, $true, 1, $false, " ²", 2)

# --- 3. Footer legend: "[1]" / "[2]" -> "¹ " / "² " -----------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    $footers = $sec.Footers
    for ($f = 1; $f -le $footers.Count; $f++) {
        $ftr = $footers.Item($f)
        if ($ftr.Exists) {
            $frng = $ftr.Range
            $frng.Find.ClearFormatting()
            $frng.Find.Replacement.ClearFormatting()
            $frng.Find.Execute("[1]", $false, $false, $false, $false, $false, $true, 1, $false, "¹ ", 2)

            $frng = $ftr.Range
            $frng.Find.ClearFormatting()
            $frng.Find.Replacement.ClearFormatting()
            $frng.Find.Execute("[2]", $false, $false, $false, $false, $false, $true, 1, $false, "² ", 2)
        }
    }
}

# --- 4. Move the _GoBack bookmark to the end of the very last edit --------
# (the final "Zwischentotal Nettovermögen insgesamt ²" occurrence, i.e. the
# position the cursor would be at after typing the last replacement).
$rng = $d.Content
$rng.Find.ClearFormatting()
$lastStart = -1
$lastEnd = -1
while ($rng.Find.Execute("Zwischentotal Nettovermögen insgesamt ²", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $rng.Start
    $lastEnd = $rng.End
    $rng.Collapse(0)
}
if ($lastStart -ge 0) {
    $target = $d.Range($lastStart, $lastEnd)
    $target.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $target)
}
